# Update plots to incorporate Courtney's suggestions:
# - Rename header "n" (B1) to "# species"
# - Rename header "count" (G1) to "# species"
# - Fill G2:G7 with the same species-count values already present in B2:B7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "# species"
$ws.Range("G1").Value = "# species"

$ws.Range("G2").Value = $ws.Range("B2").Value2
$ws.Range("G3").Value = $ws.Range("B3").Value2
$ws.Range("G4").Value = $ws.Range("B4").Value2
$ws.Range("G5").Value = $ws.Range("B5").Value2
$ws.Range("G6").Value = $ws.Range("B6").Value2
$ws.Range("G7").Value = $ws.Range("B7").Value2

$excel.Goto($ws.Range("C2:E7"))
